$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "27.206.36"
$ws.Range("E2").Value = "  +5.58%  "
Set-TextValue $ws "D3" "1.882.19"
$ws.Range("E3").Value = "  +3.90%  "
Set-TextValue $ws "D4" "1.0000"
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws "D5" "281.37"
$ws.Range("E5").Value = "  +1.63%  "
Set-TextValue $ws "D6" "0.9998"
$ws.Range("E6").Value = "  -0.05%  "
Set-TextValue $ws "D7" "0.5304"
$ws.Range("E7").Value = "  +4.06%  "
Set-TextValue $ws "D8" "0.3538"
$ws.Range("E8").Value = "  +0.46%  "
Set-TextValue $ws "D9" "45.54"
$ws.Range("E9").Value = "  +1.68%  "
Set-TextValue $ws "D10" "0.07031"
$ws.Range("E10").Value = "  +5.57%  "
Set-TextValue $ws "D11" "20.38"
$ws.Range("E11").Value = "  +1.60%  "
Set-TextValue $ws "D12" "0.8231"
$ws.Range("E12").Value = "  -1.48%  "
Set-TextValue $ws "D13" "0.07824"
$ws.Range("E13").Value = "  -0.17%  "
Set-TextValue $ws "D14" "1.883.53"
$ws.Range("E14").Value = "  +5.07%  "
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("E16").Value = "  +2.59%  "
Set-TextValue $ws "D17" "0.9998"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  +5.25%  "
$ws.Range("E19").Value = "  +2.16%  "
Set-TextValue $ws "D20" "0.9997"
$ws.Range("E20").Value = "  +0.02%  "
Set-TextValue $ws "D21" "27.232.12"
$ws.Range("E21").Value = "  +5.36%  "
Set-TextValue $ws "D22" "2.127.36"
$ws.Range("E22").Value = "  +4.75%  "
Set-TextValue $ws "D23" "4.780"
$ws.Range("E23").Value = "  +0.94%  "
Set-TextValue $ws "D24" "10.16"
$ws.Range("E24").Value = "  +1.37%  "
Set-TextValue $ws "D25" "6.258"
$ws.Range("E25").Value = "  +3.25%  "
$ws.Range("E26").Value = "  +8.96%  "
Set-TextValue $ws "D27" "147.24"
$ws.Range("E27").Value = "  +4.13%  "
Set-TextValue $ws "D28" "17.59"
Set-TextValue $ws "D29" "1.673"
$ws.Range("E29").Value = "  +1.05%  "
Set-TextValue $ws "D30" "114.93"
$ws.Range("E30").Value = "  +5.50%  "
Set-TextValue $ws "D31" "4.423"
$ws.Range("E31").Value = "  +1.74%  "
Set-TextValue $ws "D32" "4.397"
$ws.Range("E32").Value = "  +4.15%  "
Set-TextValue $ws "D33" "0.08974"
$ws.Range("E33").Value = "  +2.05%  "
Set-TextValue $ws "D34" "0.04954"
$ws.Range("E34").Value = "  +1.25%  "
Set-TextValue $ws "D35" "1.187"
$ws.Range("E35").Value = "  +4.27%  "
Set-TextValue $ws "D36" "0.7506"
$ws.Range("E36").Value = "  +2.25%  "
Set-TextValue $ws "D37" "2.909"
Set-TextValue $ws "D38" "3.316"
$ws.Range("E38").Value = "  +8.86%  "
Set-TextValue $ws "D39" "2.413"
$ws.Range("E39").Value = "  +5.55%  "
Set-TextValue $ws "D40" "0.5310"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("E41").Value = "  +1.86%  "
Set-TextValue $ws "D42" "0.9741"
$ws.Range("E42").Value = "  +2.33%  "
Set-TextValue $ws "D43" "116.99"
$ws.Range("E43").Value = "  +4.54%  "
Set-TextValue $ws "D44" "6.328"
$ws.Range("E44").Value = "  +2.30%  "
Set-TextValue $ws "D45" "8.231"
$ws.Range("E45").Value = "  +1.15%  "
Set-TextValue $ws "D46" "0.9993"
$ws.Range("E46").Value = "  -0.04%  "
Set-TextValue $ws "D47" "0.4616"
$ws.Range("E47").Value = "  +0.80%  "
Set-TextValue $ws "D48" "0.1372"
$ws.Range("E48").Value = "  -0.62%  "
Set-TextValue $ws "D49" "9.454"
$ws.Range("E49").Value = "  +1.24%  "
Set-TextValue $ws "D50" "36.74"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("E51").Value = "  +2.07%  "
